# Weekly update: a new price record (2022-02-08) is inserted as the
# second data row (row 18) of the sheet, pushing every subsequent
# record down by one row (old row 18 becomes row 19, ..., old row 82
# becomes row 83). Only the new row's values need to be written
# explicitly; Excel's row Insert shifts all the other rows (and their
# formatting) down automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 18 - shifts rows 18:82 down to 19:83
# and carries formatting down from the row above, matching the dimension
# change from A1:R82 to A1:R83.
$ws.Rows("18:18").Insert()

# Populate the newly inserted row 18 with the new record.
$ws.Range("A18").Value = 7
$ws.Range("B18").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C18").Value = "Ñuble"
$ws.Range("D18").Value = 44600
$ws.Range("E18").Value = 16
$ws.Range("F18").Value = 100112030
$ws.Range("G18").Value = "Poroto granado"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = 22000
$ws.Range("L18").Value = 23000
$ws.Range("M18").Value = 22500
$ws.Range("N18").Value = "$/saco 25 kilos"
$ws.Range("O18").Value = "Provincia de Diguillín"
$ws.Range("P18").Value = 900
$ws.Range("Q18").Value = 25
$ws.Range("R18").Value = "Hortaliza"
